$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 8 with the new bitácora entry
$ws.Range("A8").Value = "Finalización de los request de la web administradora."
$ws.Range("B8").Value = 42878.625
$ws.Range("C8").Value = 42879.041666666664
$ws.Range("D8").Value = 0.41666666666666669

# Apply the same number formats/styles as used by the existing rows above
$ws.Range("B7:D7").Copy()
$ws.Range("B8:D8").PasteSpecial(-4122)

# Update the active selection to D9, matching the diff
$ws.Range("D9").Select()
